$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.03777678762336
$ws.Cells.Item(2, 4).Value = 1.037736977939337
$ws.Cells.Item(2, 5).Value = 1.045974249498022
$ws.Cells.Item(2, 6).Value = 1.055556195354918
$ws.Cells.Item(2, 9).Value = 1.040861074483044
$ws.Cells.Item(2, 10).Value = 1.042877817621427
$ws.Cells.Item(2, 11).Value = 1.040527009585522
$ws.Cells.Item(2, 12).Value = 1.048740985483936
$ws.Cells.Item(2, 13).Value = 1.058296328623705
$ws.Cells.Item(2, 14).Value = 1.044358823742207
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.038707320684673
$ws.Cells.Item(3, 4).Value = 1.038246809317317
$ws.Cells.Item(3, 5).Value = 1.046810860112523
$ws.Cells.Item(3, 6).Value = 1.056498375303
$ws.Cells.Item(3, 9).Value = 1.041068686417052
$ws.Cells.Item(3, 10).Value = 1.043452965225869
$ws.Cells.Item(3, 11).Value = 1.040847504310623
$ws.Cells.Item(3, 12).Value = 1.049389079766284
$ws.Cells.Item(3, 13).Value = 1.059051644977546
$ws.Cells.Item(3, 14).Value = 1.044934788122218
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.039309866341272
$ws.Cells.Item(4, 4).Value = 1.03857678014702
$ws.Cells.Item(4, 5).Value = 1.047352969096884
$ws.Cells.Item(4, 6).Value = 1.057108909785256
$ws.Cells.Item(4, 9).Value = 1.041201785213006
$ws.Cells.Item(4, 10).Value = 1.043824925299288
$ws.Cells.Item(4, 11).Value = 1.041054245850161
$ws.Cells.Item(4, 12).Value = 1.04980855831473
$ws.Cells.Item(4, 13).Value = 1.059540638712862
$ws.Cells.Item(4, 14).Value = 1.045307276421606
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.039563278063453
$ws.Cells.Item(5, 4).Value = 1.038715516338208
$ws.Cells.Item(5, 5).Value = 1.047581053681491
$ws.Cells.Item(5, 6).Value = 1.057365787955308
$ws.Cells.Item(5, 9).Value = 1.041257442350817
$ws.Cells.Item(5, 10).Value = 1.043981248695186
$ws.Cells.Item(5, 11).Value = 1.041141005662485
$ws.Cells.Item(5, 12).Value = 1.049984934274584
$ws.Cells.Item(5, 13).Value = 1.059746271113388
$ws.Cells.Item(5, 14).Value = 1.045463821814656
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.039605832929426
$ws.Cells.Item(6, 4).Value = 1.038738811670487
$ws.Cells.Item(6, 5).Value = 1.047619360716596
$ws.Cells.Item(6, 6).Value = 1.057408931171817
$ws.Cells.Item(6, 9).Value = 1.0412667699566
$ws.Cells.Item(6, 10).Value = 1.04400749320096
$ws.Cells.Item(6, 11).Value = 1.041155563934111
$ws.Cells.Item(6, 12).Value = 1.050014550135466
$ws.Cells.Item(6, 13).Value = 1.059780801146735
$ws.Cells.Item(6, 14).Value = 1.045490103590638
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.039313252043499
$ws.Cells.Item(7, 4).Value = 1.038578633884293
$ws.Cells.Item(7, 5).Value = 1.047356016061027
$ws.Cells.Item(7, 6).Value = 1.057112341382671
$ws.Cells.Item(7, 9).Value = 1.041202530076676
$ws.Cells.Item(7, 10).Value = 1.043827014292456
$ws.Cells.Item(7, 11).Value = 1.041055405746845
$ws.Cells.Item(7, 12).Value = 1.049810914953555
$ws.Cells.Item(7, 13).Value = 1.059543386151116
$ws.Cells.Item(7, 14).Value = 1.045309368381383
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.038091176253729
$ws.Cells.Item(8, 4).Value = 1.037909261117153
$ws.Cells.Item(8, 5).Value = 1.046256826525584
$ws.Cells.Item(8, 6).Value = 1.055874426178102
$ws.Cells.Item(8, 9).Value = 1.04093149421785
$ws.Cells.Item(8, 10).Value = 1.04307223202327
$ws.Cells.Item(8, 11).Value = 1.040635454118475
$ws.Cells.Item(8, 12).Value = 1.048959986830144
$ws.Cells.Item(8, 13).Value = 1.058551538062029
$ws.Cells.Item(8, 14).Value = 1.044553514234801
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.035941051751986
$ws.Cells.Item(9, 4).Value = 1.03673039656713
$ws.Cells.Item(9, 5).Value = 1.044325842725719
$ws.Cells.Item(9, 6).Value = 1.053699873370194
$ws.Cells.Item(9, 9).Value = 1.040444429879395
$ws.Cells.Item(9, 10).Value = 1.041740730350977
$ws.Cells.Item(9, 11).Value = 1.039890593536804
$ws.Cells.Item(9, 12).Value = 1.047461497614457
$ws.Cells.Item(9, 13).Value = 1.056805772739615
$ws.Cells.Item(9, 14).Value = 1.043220121677407
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.034509935987934
$ws.Cells.Item(10, 4).Value = 1.035945028808505
$ws.Cells.Item(10, 5).Value = 1.043042589384148
$ws.Cells.Item(10, 6).Value = 1.05225483711058
$ws.Cells.Item(10, 9).Value = 1.040113400357296
$ws.Cells.Item(10, 10).Value = 1.040852127377096
$ws.Cells.Item(10, 11).Value = 1.039390823799981
$ws.Cells.Item(10, 12).Value = 1.046463212894964
$ws.Cells.Item(10, 13).Value = 1.055643346943075
$ws.Cells.Item(10, 14).Value = 1.042330256785379
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.033890805447029
$ws.Cells.Item(11, 4).Value = 1.035605105743912
$ws.Cells.Item(11, 5).Value = 1.042487909003223
$ws.Cells.Item(11, 6).Value = 1.051630244882418
$ws.Cells.Item(11, 9).Value = 1.039968570958407
$ws.Cells.Item(11, 10).Value = 1.040467142910711
$ws.Cells.Item(11, 11).Value = 1.039173673855599
$ws.Cells.Item(11, 12).Value = 1.046031126270599
$ws.Cells.Item(11, 13).Value = 1.055140355505609
$ws.Cells.Item(11, 14).Value = 1.041944725596893
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.033660916708025
$ws.Cells.Item(12, 4).Value = 1.035478866891984
$ws.Cells.Item(12, 5).Value = 1.042282024139466
$ws.Cells.Item(12, 6).Value = 1.051398412895544
$ws.Cells.Item(12, 9).Value = 1.03991455150139
$ws.Cells.Item(12, 10).Value = 1.040324111415499
$ws.Cells.Item(12, 11).Value = 1.039092903675648
$ws.Cells.Item(12, 12).Value = 1.04587065776856
$ws.Cells.Item(12, 13).Value = 1.054953575382379
$ws.Cells.Item(12, 14).Value = 1.041801490980552
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.033710224810093
$ws.Cells.Item(13, 4).Value = 1.035505944453796
$ws.Cells.Item(13, 5).Value = 1.04232618041649
$ws.Cells.Item(13, 6).Value = 1.051448133956425
$ws.Cells.Item(13, 9).Value = 1.039926148954989
$ws.Cells.Item(13, 10).Value = 1.04035479354979
$ws.Cells.Item(13, 11).Value = 1.03911023415964
$ws.Cells.Item(13, 12).Value = 1.045905077535746
$ws.Cells.Item(13, 13).Value = 1.054993637921679
$ws.Cells.Item(13, 14).Value = 1.041832216686993
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.033871801054743
$ws.Cells.Item(14, 4).Value = 1.035594670307442
$ws.Cells.Item(14, 5).Value = 1.042470887466457
$ws.Cells.Item(14, 6).Value = 1.05161107811468
$ws.Cells.Item(14, 9).Value = 1.039964110245858
$ws.Cells.Item(14, 10).Value = 1.040455320512068
$ws.Cells.Item(14, 11).Value = 1.039166999625743
$ws.Cells.Item(14, 12).Value = 1.046017861322663
$ws.Cells.Item(14, 13).Value = 1.055124915108155
$ws.Cells.Item(14, 14).Value = 1.041932886409089
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.033971364659968
$ws.Cells.Item(15, 4).Value = 1.035649340427263
$ws.Cells.Item(15, 5).Value = 1.042560065925285
$ws.Cells.Item(15, 6).Value = 1.051711495874618
$ws.Cells.Item(15, 9).Value = 1.039987469872356
$ws.Cells.Item(15, 10).Value = 1.040517254395831
$ws.Cells.Item(15, 11).Value = 1.039201960024521
$ws.Cells.Item(15, 12).Value = 1.046087354840119
$ws.Cells.Item(15, 13).Value = 1.055205806413672
$ws.Cells.Item(15, 14).Value = 1.04199490824607
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.034551037356582
$ws.Cells.Item(16, 4).Value = 1.035967591613419
$ws.Cells.Item(16, 5).Value = 1.043079422391416
$ws.Cells.Item(16, 6).Value = 1.052296312915069
$ws.Cells.Item(16, 9).Value = 1.040122980844601
$ws.Cells.Item(16, 10).Value = 1.04087767311122
$ws.Cells.Item(16, 11).Value = 1.039405219696589
$ws.Cells.Item(16, 12).Value = 1.046491892887705
$ws.Cells.Item(16, 13).Value = 1.055676736234791
$ws.Cells.Item(16, 14).Value = 1.042355838797376
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.034914799204998
$ws.Cells.Item(17, 4).Value = 1.036167262575405
$ws.Cells.Item(17, 5).Value = 1.043405463465765
$ws.Cells.Item(17, 6).Value = 1.052663453510451
$ws.Cells.Item(17, 9).Value = 1.040207584493621
$ws.Cells.Item(17, 10).Value = 1.041103697759415
$ws.Cells.Item(17, 11).Value = 1.039532520083319
$ws.Cells.Item(17, 12).Value = 1.046745697133221
$ws.Cells.Item(17, 13).Value = 1.055972231725715
$ws.Cells.Item(17, 14).Value = 1.042582184426496
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.035127028534027
$ws.Cells.Item(18, 4).Value = 1.036283741275722
$ws.Cells.Item(18, 5).Value = 1.043595731861491
$ws.Cells.Item(18, 6).Value = 1.05287770818526
$ws.Cells.Item(18, 9).Value = 1.040256788422202
$ws.Cells.Item(18, 10).Value = 1.041235513447493
$ws.Cells.Item(18, 11).Value = 1.039606700158981
$ws.Cells.Item(18, 12).Value = 1.046893753884932
$ws.Cells.Item(18, 13).Value = 1.056144622606574
$ws.Cells.Item(18, 14).Value = 1.042714187307971
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.035199402244307
$ws.Cells.Item(19, 4).Value = 1.03632345982572
$ws.Cells.Item(19, 5).Value = 1.043660624395612
$ws.Cells.Item(19, 6).Value = 1.052950781730551
$ws.Cells.Item(19, 9).Value = 1.040273541252815
$ws.Cells.Item(19, 10).Value = 1.041280455656482
$ws.Cells.Item(19, 11).Value = 1.039631981386956
$ws.Cells.Item(19, 12).Value = 1.046944240266769
$ws.Cells.Item(19, 13).Value = 1.056203409063391
$ws.Cells.Item(19, 14).Value = 1.042759193340052
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.034875765491758
$ws.Cells.Item(20, 4).Value = 1.036145838309676
$ws.Cells.Item(20, 5).Value = 1.043370472595441
$ws.Cells.Item(20, 6).Value = 1.052624051641068
$ws.Cells.Item(20, 9).Value = 1.040198522206962
$ws.Cells.Item(20, 10).Value = 1.041079449584315
$ws.Cells.Item(20, 11).Value = 1.039518869404779
$ws.Cells.Item(20, 12).Value = 1.046718464570416
$ws.Cells.Item(20, 13).Value = 1.055940524396632
$ws.Cells.Item(20, 14).Value = 1.042557901816207
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.033824218553684
$ws.Cells.Item(21, 4).Value = 1.035568542086698
$ws.Cells.Item(21, 5).Value = 1.042428270746894
$ws.Cells.Item(21, 6).Value = 1.051563090434409
$ws.Cells.Item(21, 9).Value = 1.039952937752326
$ws.Cells.Item(21, 10).Value = 1.040425718678718
$ws.Cells.Item(21, 11).Value = 1.03915028666668
$ws.Cells.Item(21, 12).Value = 1.045984648534905
$ws.Cells.Item(21, 13).Value = 1.055086255767812
$ws.Cells.Item(21, 14).Value = 1.04190324253774
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.033163555224761
$ws.Cells.Item(22, 4).Value = 1.035205711060178
$ws.Cells.Item(22, 5).Value = 1.041836728746337
$ws.Cells.Item(22, 6).Value = 1.050897002673444
$ws.Cells.Item(22, 9).Value = 1.039797237057166
$ws.Cells.Item(22, 10).Value = 1.040014512517445
$ws.Cells.Item(22, 11).Value = 1.038917902277681
$ws.Cells.Item(22, 12).Value = 1.045523430118303
$ws.Cells.Item(22, 13).Value = 1.054549452112247
$ws.Cells.Item(22, 14).Value = 1.041491452416552
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.033513738965031
$ws.Cells.Item(23, 4).Value = 1.035398040996693
$ws.Cells.Item(23, 5).Value = 1.042150234593892
$ws.Cells.Item(23, 6).Value = 1.051250014972547
$ws.Cells.Item(23, 9).Value = 1.03987989917432
$ws.Cells.Item(23, 10).Value = 1.040232517388131
$ws.Cells.Item(23, 11).Value = 1.039041154104286
$ws.Cells.Item(23, 12).Value = 1.04576791512015
$ws.Cells.Item(23, 13).Value = 1.054833992195738
$ws.Cells.Item(23, 14).Value = 1.041709766879159
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.034893402982276
$ws.Cells.Item(24, 4).Value = 1.0361555189701
$ws.Cells.Item(24, 5).Value = 1.043386283172799
$ws.Cells.Item(24, 6).Value = 1.052641855316932
$ws.Cells.Item(24, 9).Value = 1.040202617509068
$ws.Cells.Item(24, 10).Value = 1.041090406354828
$ws.Cells.Item(24, 11).Value = 1.039525037781688
$ws.Cells.Item(24, 12).Value = 1.046730769740819
$ws.Cells.Item(24, 13).Value = 1.055954851470855
$ws.Cells.Item(24, 14).Value = 1.042568874146591
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.036496509551003
$ws.Cells.Item(25, 4).Value = 1.037035073161799
$ws.Cells.Item(25, 5).Value = 1.04482433718635
$ws.Cells.Item(25, 6).Value = 1.054261231322308
$ws.Cells.Item(25, 9).Value = 1.040571464803526
$ws.Cells.Item(25, 10).Value = 1.042085124613451
$ws.Cells.Item(25, 11).Value = 1.040083726013622
$ws.Cells.Item(25, 12).Value = 1.047848773029724
$ws.Cells.Item(25, 13).Value = 1.057256850495874
$ws.Cells.Item(25, 14).Value = 1.043565005019237
